# RequirementsStatus.xlsx - "Add Employee Dashboard and account approval
# functionality"
#
# The "User Stories" table on Sheet1 tracks completion status per story in
# column C ("Met") and implementation notes in column D ("Notes"). This
# change marks two more stories as done (rows 14 & 15, "y" = fully met) and
# finishes the "employee can approve/reject an account registration" story
# (row 17), which moves from "not started" (no fill, empty Met/Notes) to
# "in progress" (the amber/gold fill used by rows 10/13/16) with a "~"
# status and an explanatory note.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows 14 & 15 ("view balance" / "withdraw or deposit") are now fully met.
$ws.Range("C14").Value = "y"
$ws.Range("C15").Value = "y"

# Row 17 ("employee can approve/reject account registrations") becomes
# "in progress": copy the fill/border formatting used by the other
# in-progress rows (row 13 is a same-shaped template: text/points columns
# left alone, Met+Notes columns filled in) onto row 17, then fill in its
# own Met/Notes text.
$ws.Range("A13:D13").Copy() | Out-Null
$ws.Range("A17:D17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("C17").Value = "~"
$ws.Range("D17").Value = "Works, but can't print username for some reason???"

# Leave the selection on the cell that was just edited.
$ws.Range("C16").Select() | Out-Null
